# Generate Report for Handback
#
# Refreshes the handback-status report: a new handoff/handback round-trip
# was recorded for the "3e55b075-..." source file (in both the zh-cn and
# de-de language sheets), so the "Correspond Handoff Datetime" and
# "Correspond Handback DateTime" columns for that row move forward, and
# the Overview sheet's "Latest HO Xliff Generate Date" (driven off the
# de-de handoff datetime) is refreshed to match.
#
# The "cee6e158-..." row is untouched in this run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-28 12:48:59"
$zhcn.Range("K2").Value = "2016-08-28 12:49:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-28 12:49:07"
$dede.Range("K2").Value = "2016-08-28 12:49:24"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-28 12:49:07"
